$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 127, shifting the existing
# rows 127-166 down to 128-167 (matches the target dimension A1:R167).
$ws.Rows(127).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A127").Value = 7
$ws.Range("B127").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C127").Value = "Ñuble"
$ws.Range("D127").Value = 44559
$ws.Range("E127").Value = 16
$ws.Range("F127").Value = 100112032
$ws.Range("G127").Value = "Zapallo italiano"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 100
$ws.Range("K127").Value = 5000
$ws.Range("L127").Value = 5500
$ws.Range("M127").Value = 5250
$ws.Range("N127").Value = "$/caja 60 unidades"
$ws.Range("O127").Value = "Región del Maule"
$ws.Range("P127").Value = 88
$ws.Range("Q127").Value = 60
$ws.Range("R127").Value = "Hortaliza"
